$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first 4 data rows (rows 2-5), shifting all subsequent rows up.
$ws.Range("A2:E5").EntireRow.Delete()
